$d = $word.ActiveDocument

# The second paragraph holds a Word field "{ m:self.name }" expressed via
# fldChar begin / instrText runs / fldChar end. The commit converts it into
# plain literal-text runs spelling "{m:self.name}" (no surrounding spaces),
# keeping the existing orange color formatting on the "self" run.

$f = $d.Fields.Item(1)
$codeStart = $f.Code.Start

# Locate the paragraph that contains the field so we know where to rebuild
# the literal-text runs after removing the field.
$target_index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($codeStart -ge $pp.Range.Start -and $codeStart -lt $pp.Range.End) {
        $target_index = $i
    }
}

# Remove the field entirely (fldChar begin/end + instrText runs).
$f.Delete()

# Re-fetch the (now field-less) paragraph; its Range now only spans the
# paragraph mark.
$p = $d.Paragraphs.Item($target_index)

# Range covering just the paragraph's content (excluding the paragraph
# mark), so inserting OOXML here only replaces the run content and keeps
# the paragraph's own properties/attributes untouched.
$target = $d.Range($p.Range.Start, $p.Range.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>self</w:t></w:r><w:r><w:t xml:space="preserve">.name}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
